# Applies the "巢湖·原铁崩ONLY" insertion + follow-on count bumps to both the
# "展览" sheet and the "全部类型" roll-up sheet (which mirrors "展览" plus the
# single row coming from "演出"). Both sheets get an identical transform for
# rows 1-21; "全部类型" simply has one extra unaffected row tacked on, and the
# row-insert naturally cascades it from row 21 to row 22.

$wb = $excel.ActiveWorkbook

function Update-ExpoSheet {
    param($ws)

    # --- small "want-to-go" counter bumps on existing rows (rows 1-8 block) ---
    $ws.Range("F4").Value = 868
    $ws.Range("F5").Value = 35
    $ws.Range("F7").Value = 10431

    # --- insert a brand-new event row at row 9, pushing rows 9.. down by one ---
    $ws.Rows.Item(9).Insert()

    # Copy formatting (border/bold/alignment on col A, etc.) down from the row
    # that is now directly above the blank inserted row, then stamp the new
    # row's values over the top.
    $ws.Range("A8:I8").Copy($ws.Range("A9:I9"))

    # Leading apostrophe forces text so Excel doesn't silently reinterpret the
    # ISO-ish "2024-05-01" literal as a date serial (the sheet stores these as
    # plain strings, matching every other date cell in column B).
    $ws.Range("B9").Value = "'2024-05-01"
    $ws.Range("C9").Value = "巢湖·原铁崩ONLY"
    $ws.Range("D9").Value = "团结东路7号 巢湖宾馆"
    $ws.Range("E9").Value = "2024.05.01 10:00-05.01 17:00"
    $ws.Range("F9").Value = 1
    $ws.Range("G9").Value = 55
    $ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=84289"
    $ws.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202404/ujpuWAig1713161918045.jpeg"

    # --- "want-to-go" counter bumps on rows that shifted down by one ---
    $ws.Range("F12").Value = 134   # was row 11 (F=131) before the insert
    $ws.Range("F13").Value = 13    # was row 12 (F=12)
    $ws.Range("F17").Value = 30    # was row 16 (F=29)
    $ws.Range("F18").Value = 289   # was row 17 (F=287)
    $ws.Range("F19").Value = 884   # was row 18 (F=830)
    $ws.Range("F20").Value = 49    # was row 19 (F=48)
    $ws.Range("F21").Value = 100   # was row 20 (F=99)

    # --- re-number the running index in column A (row number - 1) for every
    #     row from the inserted one through the sheet's last used row ---
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 9; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

Update-ExpoSheet $wb.Worksheets.Item("展览")
Update-ExpoSheet $wb.Worksheets.Item("全部类型")
